$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Short Term" ---
$ws1 = $wb.Worksheets.Item("Short Term")

# Row 113
$ws1.Range("B113").Value = 1.16
$ws1.Range("C113").Value = 5.27
$ws1.Range("D113").Value = -2.74
$ws1.Range("E113").Value = 25.86
$ws1.Range("F113").Value = 27.97
$ws1.Range("G113").Value = -0.34

# Row 114
$ws1.Range("B114").Value = -2.72
$ws1.Range("C114").Value = -9.85
$ws1.Range("D114").Value = -1.77

# Row 115
$ws1.Range("B115").Value = -2.45
$ws1.Range("C115").Value = -0.85
$ws1.Range("D115").Value = 6.35

# Row 116
$ws1.Range("B116").Value = 5.91
$ws1.Range("C116").Value = 3.36
$ws1.Range("D116").Value = 3.34

# Row 117
$ws1.Range("B117").Value = 0.84
$ws1.Range("C117").Value = 1.91
$ws1.Range("D117").Value = 6.35

# Row 118
$ws1.Range("B118").Value = 19.76
$ws1.Range("C118").Value = 23.75
$ws1.Range("D118").Value = -11.65
$ws1.Range("E118").Value = 41.36
$ws1.Range("F118").Value = 41.83
$ws1.Range("G118").Value = -4.7

# Row 119 (new row)
$ws1.Range("A118").Copy($ws1.Range("A119"))
$ws1.Range("A119").Value = 45566
$ws1.Range("B119").Value = -9.82
$ws1.Range("C119").Value = -7.06
$ws1.Range("D119").Value = -2.95
$ws1.Range("E119").Value = 24.9
$ws1.Range("F119").Value = 20.9
$ws1.Range("G119").Value = -1.85

# --- Sheet 2: "Medium Term" ---
$ws2 = $wb.Worksheets.Item("Medium Term")

# Row 99
$ws2.Range("B99").Value = 7.7
$ws2.Range("C99").Value = 4.79
$ws2.Range("D99").Value = -1.03

# Row 100
$ws2.Range("B100").Value = 12.07
$ws2.Range("C100").Value = 7.34
$ws2.Range("D100").Value = 0.3

# Row 101
$ws2.Range("B101").Value = 8.1
$ws2.Range("C101").Value = 7.05
$ws2.Range("D101").Value = -0.81

# Row 102
$ws2.Range("C102").Value = 5.92
$ws2.Range("D102").Value = -0.05

# Row 103
$ws2.Range("C103").Value = 7.92
$ws2.Range("D103").Value = 2.17

# Row 104
$ws2.Range("B104").Value = 21.23
$ws2.Range("C104").Value = 14.49
$ws2.Range("D104").Value = 7.35

# Row 105 (new row)
$ws2.Range("A104").Copy($ws2.Range("A105"))
$ws2.Range("A105").Value = 45566
$ws2.Range("B105").Value = 23.22
$ws2.Range("C105").Value = 13.57
$ws2.Range("D105").Value = 9.21
